# "finish import from excel"
# The sheet used to hold a 3-user export (Id/Name/Email/password/Age/Created At/
# Updated At headers + noor/naser/salim rows). The finished import collapses
# this down to a single compact row: Name, Email (as a mailto hyperlink),
# password and Age - with no header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old per-row mailto hyperlinks (C2:C4) before we touch the data so
# we don't drag stale hyperlink formatting/relationships along.
$ws.Hyperlinks.Delete()

# Wipe the old 7-column x 4-row table entirely.
$ws.Cells.Clear()

# Write the new, single-row data set.
$ws.Range("A1").Value = "Noor"
$ws.Range("B1").Value = "noor@noor.com"
$ws.Range("C1").Value = "Noo12345678"
$ws.Range("D1").Value = 20

# Re-create the hyperlink on the email cell and restore the standard
# "Hyperlink" cell style (underline + theme color) on it.
$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:noor@noor.com")
$ws.Range("B1").Style = "Hyperlink"

# Match the saved selection/active cell from the finished workbook.
$ws.Range("E3").Select() | Out-Null
